$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 562, shifting the existing rows
# (562-623) down to (563-624).
$ws.Rows.Item(562).Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Range("A562").Value = 3
$ws.Range("B562").Value = "Femacal de La Calera"
$ws.Range("C562").Value = "Coquimbo"
$ws.Range("D562").Value = 45194
$ws.Range("E562").Value = 5
$ws.Range("F562").Value = 100112009
$ws.Range("G562").Value = "Acelga"
$ws.Range("H562").Value = "Sin especificar"
$ws.Range("I562").Value = "Primera"
$ws.Range("J562").Value = 150
$ws.Range("K562").Value = 3500
$ws.Range("L562").Value = 3800
$ws.Range("M562").Value = 3580
$ws.Range("N562").Value = "`$/docena de atados (6 kilos)"
$ws.Range("O562").Value = "Provincia de Quillota"
$ws.Range("P562").Value = 597
$ws.Range("Q562").Value = 6
$ws.Range("R562").Value = "Hortaliza"
